$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Calr"
$ws.Cells.Item(2,3).Value = "Scarf1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 82.73729466666667
$ws.Cells.Item(2,8).Value = 248.211884
$ws.Cells.Item(2,9).Value = 0.09847102321391109
$ws.Cells.Item(2,10).Value = 0.09847102321391106
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 15.16337366666667
$ws.Cells.Item(2,14).Value = 45.490121
$ws.Cells.Item(2,15).Value = 0.722233085414401
$ws.Cells.Item(2,16).Value = 0.7222330854144011
$ws.Cells.Item(2,17).Value = 1254.576515199774
$ws.Cells.Item(2,18).Value = 11291.18863679796
$ws.Cells.Item(2,19).Value = 0.07111903091969611
$ws.Cells.Item(2,20).Value = 0.0711190309196961
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Calr"
$ws.Cells.Item(3,3).Value = "Scarf1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 82.73729466666667
$ws.Cells.Item(3,8).Value = 248.211884
$ws.Cells.Item(3,9).Value = 0.09847102321391109
$ws.Cells.Item(3,10).Value = 0.09847102321391106
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.2227456666666666
$ws.Cells.Item(3,14).Value = 0.668237
$ws.Cells.Item(3,15).Value = 0.01060939957266904
$ws.Cells.Item(3,16).Value = 0.01060939957266904
$ws.Cells.Item(3,17).Value = 18.42937385872311
$ws.Cells.Item(3,18).Value = 165.864364728508
$ws.Cells.Item(3,19).Value = 0.001044718431605952
$ws.Cells.Item(3,20).Value = 0.001044718431605952
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Calr"
$ws.Cells.Item(4,3).Value = "Scarf1"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 82.73729466666667
$ws.Cells.Item(4,8).Value = 248.211884
$ws.Cells.Item(4,9).Value = 0.09847102321391109
$ws.Cells.Item(4,10).Value = 0.09847102321391106
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.016357666666666
$ws.Cells.Item(4,14).Value = 6.049073
$ws.Cells.Item(4,15).Value = 0.09603932811449209
$ws.Cells.Item(4,16).Value = 0.0960393281144921
$ws.Cells.Item(4,17).Value = 166.8279784203924
$ws.Cells.Item(4,18).Value = 1501.451805783532
$ws.Cells.Item(4,19).Value = 0.009457090908210574
$ws.Cells.Item(4,20).Value = 0.009457090908210572
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Calr"
$ws.Cells.Item(5,3).Value = "Scarf1"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 82.73729466666667
$ws.Cells.Item(5,8).Value = 248.211884
$ws.Cells.Item(5,9).Value = 0.09847102321391109
$ws.Cells.Item(5,10).Value = 0.09847102321391106
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.163758666666666
$ws.Cells.Item(5,14).Value = 9.491275999999999
$ws.Cells.Item(5,15).Value = 0.1506901586390516
$ws.Cells.Item(5,16).Value = 0.1506901586390516
$ws.Cells.Item(5,17).Value = 261.7608330582204
$ws.Cells.Item(5,18).Value = 2355.847497523984
$ws.Cells.Item(5,19).Value = 0.014838614109454
$ws.Cells.Item(5,20).Value = 0.01483861410945399
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Calr"
$ws.Cells.Item(6,3).Value = "Scarf1"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 82.73729466666667
$ws.Cells.Item(6,8).Value = 248.211884
$ws.Cells.Item(6,9).Value = 0.09847102321391109
$ws.Cells.Item(6,10).Value = 0.09847102321391106
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.428889
$ws.Cells.Item(6,14).Value = 1.286667
$ws.Cells.Item(6,15).Value = 0.02042802825938606
$ws.Cells.Item(6,16).Value = 0.02042802825938606
$ws.Cells.Item(6,17).Value = 35.485115572292
$ws.Cells.Item(6,18).Value = 319.366040150628
$ws.Cells.Item(6,19).Value = 0.002011568844944437
$ws.Cells.Item(6,20).Value = 0.002011568844944436
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Calr"
$ws.Cells.Item(7,3).Value = "Scarf1"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 237.0718893333334
$ws.Cells.Item(7,8).Value = 711.2156680000001
$ws.Cells.Item(7,9).Value = 0.2821546391135941
$ws.Cells.Item(7,10).Value = 0.2821546391135941
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 15.16337366666667
$ws.Cells.Item(7,14).Value = 45.490121
$ws.Cells.Item(7,15).Value = 0.722233085414401
$ws.Cells.Item(7,16).Value = 0.7222330854144011
$ws.Cells.Item(7,17).Value = 3594.809643823981
$ws.Cells.Item(7,18).Value = 32353.28679441583
$ws.Cells.Item(7,19).Value = 0.2037814155709979
$ws.Cells.Item(7,20).Value = 0.2037814155709979
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Calr"
$ws.Cells.Item(8,3).Value = "Scarf1"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 237.0718893333334
$ws.Cells.Item(8,8).Value = 711.2156680000001
$ws.Cells.Item(8,9).Value = 0.2821546391135941
$ws.Cells.Item(8,10).Value = 0.2821546391135941
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.2227456666666666
$ws.Cells.Item(8,14).Value = 0.668237
$ws.Cells.Item(8,15).Value = 0.01060939957266904
$ws.Cells.Item(8,16).Value = 0.01060939957266904
$ws.Cells.Item(8,17).Value = 52.80673603747956
$ws.Cells.Item(8,18).Value = 475.260624337316
$ws.Cells.Item(8,19).Value = 0.002993491307638353
$ws.Cells.Item(8,20).Value = 0.002993491307638353
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Calr"
$ws.Cells.Item(9,3).Value = "Scarf1"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 237.0718893333334
$ws.Cells.Item(9,8).Value = 711.2156680000001
$ws.Cells.Item(9,9).Value = 0.2821546391135941
$ws.Cells.Item(9,10).Value = 0.2821546391135941
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.016357666666666
$ws.Cells.Item(9,14).Value = 6.049073
$ws.Cells.Item(9,15).Value = 0.09603932811449209
$ws.Cells.Item(9,16).Value = 0.0960393281144921
$ws.Cells.Item(9,17).Value = 478.0217216084183
$ws.Cells.Item(9,18).Value = 4302.195494475764
$ws.Cells.Item(9,19).Value = 0.02709794196485657
$ws.Cells.Item(9,20).Value = 0.02709794196485657
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Calr"
$ws.Cells.Item(10,3).Value = "Scarf1"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 237.0718893333334
$ws.Cells.Item(10,8).Value = 711.2156680000001
$ws.Cells.Item(10,9).Value = 0.2821546391135941
$ws.Cells.Item(10,10).Value = 0.2821546391135941
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.163758666666666
$ws.Cells.Item(10,14).Value = 9.491275999999999
$ws.Cells.Item(10,15).Value = 0.1506901586390516
$ws.Cells.Item(10,16).Value = 0.1506901586390516
$ws.Cells.Item(10,17).Value = 750.0382445013743
$ws.Cells.Item(10,18).Value = 6750.344200512368
$ws.Cells.Item(10,19).Value = 0.04251792732877186
$ws.Cells.Item(10,20).Value = 0.04251792732877185
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Calr"
$ws.Cells.Item(11,3).Value = "Scarf1"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 237.0718893333334
$ws.Cells.Item(11,8).Value = 711.2156680000001
$ws.Cells.Item(11,9).Value = 0.2821546391135941
$ws.Cells.Item(11,10).Value = 0.2821546391135941
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.428889
$ws.Cells.Item(11,14).Value = 1.286667
$ws.Cells.Item(11,15).Value = 0.02042802825938606
$ws.Cells.Item(11,16).Value = 0.02042802825938606
$ws.Cells.Item(11,17).Value = 101.677525544284
$ws.Cells.Item(11,18).Value = 915.0977298985561
$ws.Cells.Item(11,19).Value = 0.005763862941329375
$ws.Cells.Item(11,20).Value = 0.005763862941329374
$ws.Cells.Item(12,1).Value = "M1"
$ws.Cells.Item(12,2).Value = "Calr"
$ws.Cells.Item(12,3).Value = "Scarf1"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 277.340215
$ws.Cells.Item(12,8).Value = 832.0206450000001
$ws.Cells.Item(12,9).Value = 0.3300805865050695
$ws.Cells.Item(12,10).Value = 0.3300805865050694
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 15.16337366666667
$ws.Cells.Item(12,14).Value = 45.490121
$ws.Cells.Item(12,15).Value = 0.722233085414401
$ws.Cells.Item(12,16).Value = 0.7222330854144011
$ws.Cells.Item(12,17).Value = 4205.413312838672
$ws.Cells.Item(12,18).Value = 37848.71981554805
$ws.Cells.Item(12,19).Value = 0.2383951204269514
$ws.Cells.Item(12,20).Value = 0.2383951204269514
$ws.Cells.Item(13,1).Value = "M1"
$ws.Cells.Item(13,2).Value = "Calr"
$ws.Cells.Item(13,3).Value = "Scarf1"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 277.340215
$ws.Cells.Item(13,8).Value = 832.0206450000001
$ws.Cells.Item(13,9).Value = 0.3300805865050695
$ws.Cells.Item(13,10).Value = 0.3300805865050694
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.2227456666666666
$ws.Cells.Item(13,14).Value = 0.668237
$ws.Cells.Item(13,15).Value = 0.01060939957266904
$ws.Cells.Item(13,16).Value = 0.01060939957266904
$ws.Cells.Item(13,17).Value = 61.77633108365166
$ws.Cells.Item(13,18).Value = 555.9869797528651
$ws.Cells.Item(13,19).Value = 0.003501956833413231
$ws.Cells.Item(13,20).Value = 0.003501956833413231
$ws.Cells.Item(14,1).Value = "M1"
$ws.Cells.Item(14,2).Value = "Calr"
$ws.Cells.Item(14,3).Value = "Scarf1"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 277.340215
$ws.Cells.Item(14,8).Value = 832.0206450000001
$ws.Cells.Item(14,9).Value = 0.3300805865050695
$ws.Cells.Item(14,10).Value = 0.3300805865050694
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 2.016357666666666
$ws.Cells.Item(14,14).Value = 6.049073
$ws.Cells.Item(14,15).Value = 0.09603932811449209
$ws.Cells.Item(14,16).Value = 0.0960393281144921
$ws.Cells.Item(14,17).Value = 559.2170687902316
$ws.Cells.Item(14,18).Value = 5032.953619112085
$ws.Cells.Item(14,19).Value = 0.03170071775158436
$ws.Cells.Item(14,20).Value = 0.03170071775158436
$ws.Cells.Item(15,1).Value = "M1"
$ws.Cells.Item(15,2).Value = "Calr"
$ws.Cells.Item(15,3).Value = "Scarf1"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 277.340215
$ws.Cells.Item(15,8).Value = 832.0206450000001
$ws.Cells.Item(15,9).Value = 0.3300805865050695
$ws.Cells.Item(15,10).Value = 0.3300805865050694
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 3.163758666666666
$ws.Cells.Item(15,14).Value = 9.491275999999999
$ws.Cells.Item(15,15).Value = 0.1506901586390516
$ws.Cells.Item(15,16).Value = 0.1506901586390516
$ws.Cells.Item(15,17).Value = 877.4375088214466
$ws.Cells.Item(15,18).Value = 7896.937579393019
$ws.Cells.Item(15,19).Value = 0.04973989594412013
$ws.Cells.Item(15,20).Value = 0.04973989594412012
$ws.Cells.Item(16,1).Value = "M1"
$ws.Cells.Item(16,2).Value = "Calr"
$ws.Cells.Item(16,3).Value = "Scarf1"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 277.340215
$ws.Cells.Item(16,8).Value = 832.0206450000001
$ws.Cells.Item(16,9).Value = 0.3300805865050695
$ws.Cells.Item(16,10).Value = 0.3300805865050694
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.428889
$ws.Cells.Item(16,14).Value = 1.286667
$ws.Cells.Item(16,15).Value = 0.02042802825938606
$ws.Cells.Item(16,16).Value = 0.02042802825938606
$ws.Cells.Item(16,17).Value = 118.948167471135
$ws.Cells.Item(16,18).Value = 1070.533507240215
$ws.Cells.Item(16,19).Value = 0.006742895549000284
$ws.Cells.Item(16,20).Value = 0.006742895549000283
$ws.Cells.Item(17,1).Value = "M2"
$ws.Cells.Item(17,2).Value = "Calr"
$ws.Cells.Item(17,3).Value = "Scarf1"
$ws.Cells.Item(17,4).Value = "ECs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 191.058024
$ws.Cells.Item(17,8).Value = 573.174072
$ws.Cells.Item(17,9).Value = 0.22739055213619
$ws.Cells.Item(17,10).Value = 0.22739055213619
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 15.16337366666667
$ws.Cells.Item(17,14).Value = 45.490121
$ws.Cells.Item(17,15).Value = 0.722233085414401
$ws.Cells.Item(17,16).Value = 0.7222330854144011
$ws.Cells.Item(17,17).Value = 2897.084209926968
$ws.Cells.Item(17,18).Value = 26073.75788934271
$ws.Cells.Item(17,19).Value = 0.1642289800634047
$ws.Cells.Item(17,20).Value = 0.1642289800634047
$ws.Cells.Item(18,1).Value = "M2"
$ws.Cells.Item(18,2).Value = "Calr"
$ws.Cells.Item(18,3).Value = "Scarf1"
$ws.Cells.Item(18,4).Value = "FAPs"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 191.058024
$ws.Cells.Item(18,8).Value = 573.174072
$ws.Cells.Item(18,9).Value = 0.22739055213619
$ws.Cells.Item(18,10).Value = 0.22739055213619
$ws.Cells.Item(18,11).Value = 2
$ws.Cells.Item(18,12).Value = 0.6666666666666666
$ws.Cells.Item(18,13).Value = 0.2227456666666666
$ws.Cells.Item(18,14).Value = 0.668237
$ws.Cells.Item(18,15).Value = 0.01060939957266904
$ws.Cells.Item(18,16).Value = 0.01060939957266904
$ws.Cells.Item(18,17).Value = 42.557346927896
$ws.Cells.Item(18,18).Value = 383.016122351064
$ws.Cells.Item(18,19).Value = 0.002412477226662672
$ws.Cells.Item(18,20).Value = 0.002412477226662672
$ws.Cells.Item(19,1).Value = "M2"
$ws.Cells.Item(19,2).Value = "Calr"
$ws.Cells.Item(19,3).Value = "Scarf1"
$ws.Cells.Item(19,4).Value = "M1"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 191.058024
$ws.Cells.Item(19,8).Value = 573.174072
$ws.Cells.Item(19,9).Value = 0.22739055213619
$ws.Cells.Item(19,10).Value = 0.22739055213619
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 2.016357666666666
$ws.Cells.Item(19,14).Value = 6.049073
$ws.Cells.Item(19,15).Value = 0.09603932811449209
$ws.Cells.Item(19,16).Value = 0.0960393281144921
$ws.Cells.Item(19,17).Value = 385.241311470584
$ws.Cells.Item(19,18).Value = 3467.171803235256
$ws.Cells.Item(19,19).Value = 0.02183843584674307
$ws.Cells.Item(19,20).Value = 0.02183843584674307
$ws.Cells.Item(20,1).Value = "M2"
$ws.Cells.Item(20,2).Value = "Calr"
$ws.Cells.Item(20,3).Value = "Scarf1"
$ws.Cells.Item(20,4).Value = "M2"
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 191.058024
$ws.Cells.Item(20,8).Value = 573.174072
$ws.Cells.Item(20,9).Value = 0.22739055213619
$ws.Cells.Item(20,10).Value = 0.22739055213619
$ws.Cells.Item(20,11).Value = 3
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 3.163758666666666
$ws.Cells.Item(20,14).Value = 9.491275999999999
$ws.Cells.Item(20,15).Value = 0.1506901586390516
$ws.Cells.Item(20,16).Value = 0.1506901586390516
$ws.Cells.Item(20,17).Value = 604.461479266208
$ws.Cells.Item(20,18).Value = 5440.153313395872
$ws.Cells.Item(20,19).Value = 0.03426551837442401
$ws.Cells.Item(20,20).Value = 0.03426551837442401
$ws.Cells.Item(21,1).Value = "M2"
$ws.Cells.Item(21,2).Value = "Calr"
$ws.Cells.Item(21,3).Value = "Scarf1"
$ws.Cells.Item(21,4).Value = "sCs"
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 191.058024
$ws.Cells.Item(21,8).Value = 573.174072
$ws.Cells.Item(21,9).Value = 0.22739055213619
$ws.Cells.Item(21,10).Value = 0.22739055213619
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 0.428889
$ws.Cells.Item(21,14).Value = 1.286667
$ws.Cells.Item(21,15).Value = 0.02042802825938606
$ws.Cells.Item(21,16).Value = 0.02042802825938606
$ws.Cells.Item(21,17).Value = 81.94268485533601
$ws.Cells.Item(21,18).Value = 737.4841636980241
$ws.Cells.Item(21,19).Value = 0.004645140624955488
$ws.Cells.Item(21,20).Value = 0.004645140624955487
$ws.Cells.Item(22,1).Value = "sCs"
$ws.Cells.Item(22,2).Value = "Calr"
$ws.Cells.Item(22,3).Value = "Scarf1"
$ws.Cells.Item(22,4).Value = "ECs"
$ws.Cells.Item(22,5).Value = 3
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,7).Value = 52.01228800000001
$ws.Cells.Item(22,8).Value = 156.036864
$ws.Cells.Item(22,9).Value = 0.0619031990312353
$ws.Cells.Item(22,10).Value = 0.06190319903123528
$ws.Cells.Item(22,11).Value = 3
$ws.Cells.Item(22,12).Value = 1
$ws.Cells.Item(22,13).Value = 15.16337366666667
$ws.Cells.Item(22,14).Value = 45.490121
$ws.Cells.Item(22,15).Value = 0.722233085414401
$ws.Cells.Item(22,16).Value = 0.7222330854144011
$ws.Cells.Item(22,17).Value = 788.6817582022827
$ws.Cells.Item(22,18).Value = 7098.135823820545
$ws.Cells.Item(22,19).Value = 0.04470853843335083
$ws.Cells.Item(22,20).Value = 0.04470853843335083
$ws.Cells.Item(23,1).Value = "sCs"
$ws.Cells.Item(23,2).Value = "Calr"
$ws.Cells.Item(23,3).Value = "Scarf1"
$ws.Cells.Item(23,4).Value = "FAPs"
$ws.Cells.Item(23,5).Value = 3
$ws.Cells.Item(23,6).Value = 1
$ws.Cells.Item(23,7).Value = 52.01228800000001
$ws.Cells.Item(23,8).Value = 156.036864
$ws.Cells.Item(23,9).Value = 0.0619031990312353
$ws.Cells.Item(23,10).Value = 0.06190319903123528
$ws.Cells.Item(23,11).Value = 2
$ws.Cells.Item(23,12).Value = 0.6666666666666666
$ws.Cells.Item(23,13).Value = 0.2227456666666666
$ws.Cells.Item(23,14).Value = 0.668237
$ws.Cells.Item(23,15).Value = 0.01060939957266904
$ws.Cells.Item(23,16).Value = 0.01060939957266904
$ws.Cells.Item(23,17).Value = 11.58551176541867
$ws.Cells.Item(23,18).Value = 104.269605888768
$ws.Cells.Item(23,19).Value = 0.0006567557733488345
$ws.Cells.Item(23,20).Value = 0.0006567557733488344
$ws.Cells.Item(24,1).Value = "sCs"
$ws.Cells.Item(24,2).Value = "Calr"
$ws.Cells.Item(24,3).Value = "Scarf1"
$ws.Cells.Item(24,4).Value = "M1"
$ws.Cells.Item(24,5).Value = 3
$ws.Cells.Item(24,6).Value = 1
$ws.Cells.Item(24,7).Value = 52.01228800000001
$ws.Cells.Item(24,8).Value = 156.036864
$ws.Cells.Item(24,9).Value = 0.0619031990312353
$ws.Cells.Item(24,10).Value = 0.06190319903123528
$ws.Cells.Item(24,11).Value = 3
$ws.Cells.Item(24,12).Value = 1
$ws.Cells.Item(24,13).Value = 2.016357666666666
$ws.Cells.Item(24,14).Value = 6.049073
$ws.Cells.Item(24,15).Value = 0.09603932811449209
$ws.Cells.Item(24,16).Value = 0.0960393281144921
$ws.Cells.Item(24,17).Value = 104.8753756696747
$ws.Cells.Item(24,18).Value = 943.878381027072
$ws.Cells.Item(24,19).Value = 0.005945141643097515
$ws.Cells.Item(24,20).Value = 0.005945141643097515
$ws.Cells.Item(25,1).Value = "sCs"
$ws.Cells.Item(25,2).Value = "Calr"
$ws.Cells.Item(25,3).Value = "Scarf1"
$ws.Cells.Item(25,4).Value = "M2"
$ws.Cells.Item(25,5).Value = 3
$ws.Cells.Item(25,6).Value = 1
$ws.Cells.Item(25,7).Value = 52.01228800000001
$ws.Cells.Item(25,8).Value = 156.036864
$ws.Cells.Item(25,9).Value = 0.0619031990312353
$ws.Cells.Item(25,10).Value = 0.06190319903123528
$ws.Cells.Item(25,11).Value = 3
$ws.Cells.Item(25,12).Value = 1
$ws.Cells.Item(25,13).Value = 3.163758666666666
$ws.Cells.Item(25,14).Value = 9.491275999999999
$ws.Cells.Item(25,15).Value = 0.1506901586390516
$ws.Cells.Item(25,16).Value = 0.1506901586390516
$ws.Cells.Item(25,17).Value = 164.5543269331627
$ws.Cells.Item(25,18).Value = 1480.988942398464
$ws.Cells.Item(25,19).Value = 0.009328202882281634
$ws.Cells.Item(25,20).Value = 0.009328202882281632
$ws.Cells.Item(26,1).Value = "sCs"
$ws.Cells.Item(26,2).Value = "Calr"
$ws.Cells.Item(26,3).Value = "Scarf1"
$ws.Cells.Item(26,4).Value = "sCs"
$ws.Cells.Item(26,5).Value = 3
$ws.Cells.Item(26,6).Value = 1
$ws.Cells.Item(26,7).Value = 52.01228800000001
$ws.Cells.Item(26,8).Value = 156.036864
$ws.Cells.Item(26,9).Value = 0.0619031990312353
$ws.Cells.Item(26,10).Value = 0.06190319903123528
$ws.Cells.Item(26,11).Value = 3
$ws.Cells.Item(26,12).Value = 1
$ws.Cells.Item(26,13).Value = 0.428889
$ws.Cells.Item(26,14).Value = 1.286667
$ws.Cells.Item(26,15).Value = 0.02042802825938606
$ws.Cells.Item(26,16).Value = 0.02042802825938606
$ws.Cells.Item(26,17).Value = 22.307498188032
$ws.Cells.Item(26,18).Value = 200.767483692288
$ws.Cells.Item(26,19).Value = 0.001264560299156474
$ws.Cells.Item(26,20).Value = 0.001264560299156474
